{"js": "// Office.js (Word JavaScript API) edit script.\n// Rewrites the \"Project Description\" body text for the \"Dash Warning\"\n// project-expo table: the sentence about scanning dashboard symbols /\n// GPS-based repair shops is replaced with the new wording about\n// learning dashboard fault symbols, and the Q&A Forum sentence is\n// merged with the roadside-assistance sentence that used to precede it.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in the document.\");\n}\n\n// The description lives in the merged cell at row index 5, column index 1\n// (0-based) of the first table - the \"Project Description\" row.\nconst table = tables.items[0];\nconst cell = table.getCell(5, 1);\nconst paragraph = cell.body.paragraphs.getFirst();\n\nconst oldText =\n  \"allows users to scan the symbols that appear on their vehicle's \" +\n  \"dashboard and displays what each symbol means. It also shows or \" +\n  \"recommends nearby repair shops in situations when a mechanic is \" +\n  \"needed using GPS location. The app will also recommend phone \" +\n  \"numbers to various roadside assistance personnel as the users \" +\n  \"request when they have broken down. The app also allows users to \" +\n  \"manually select the type of fault they may be experiencing by \" +\n  \"using a Q&A Forum.\";\n\nconst newText =\n  \"allows users to learn about the fault symbols that appear on their \" +\n  \"vehicle's dashboard, by displaying what the fault they've \" +\n  \"experienced means. The app allows users to manually select the \" +\n  \"type of fault they may be experiencing using a Q&A Forum. The app \" +\n  \"will also recommend phone numbers to various roadside assistance \" +\n  \"personnel, as the users request, when they have broken down. \";\n\nconst searchResults = paragraph.search(oldText, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length !== 1) {\n  throw new Error(\n    \"Expected to find exactly one occurrence of the target sentence, found \" +\n      searchResults.items.length\n  );\n}\n\nsearchResults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Rewrites the \"Project Description\" body text for the \"Dash Warning\"\n# project-expo table: the sentence about scanning dashboard symbols /\n# GPS-based repair shops is replaced with the new wording about\n# learning dashboard fault symbols, and the Q&A Forum sentence is\n# merged with the roadside-assistance sentence that used to precede it.\n\n$d = $word.ActiveDocument\n\n# Disable smart-quote autocorrect so the straight apostrophes in the\n# replacement text are not turned into curly quotes on save.\n$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false\n$word.Options.AutoFormatReplaceQuotes = $false\n\n$oldText = \"allows users to scan the symbols that appear on their vehicle's dashboard and displays what each symbol means. It also shows or recommends nearby repair shops in situations when a mechanic is needed using GPS location. The app will also recommend phone numbers to various roadside assistance personnel as the users request when they have broken down. The app also allows users to manually select the type of fault they may be experiencing by using a Q&A Forum.\"\n$newText = \"allows users to learn about the fault symbols that appear on their vehicle's dashboard, by displaying what the fault they've experienced means. The app allows users to manually select the type of fault they may be experiencing using a Q&A Forum. The app will also recommend phone numbers to various roadside assistance personnel, as the users request, when they have broken down. \"\n\n$searchRange = $d.Content\n$searchRange.Find.Text = $oldText\n$found = $searchRange.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the target sentence to replace.\"\n}\n\n# Re-anchor a plain document Range over the found span and assign its\n# .Text directly - this replaces the span in place without disturbing\n# the rest of the paragraph (Find/Execute's own Replacement.Text path\n# is what mangles the apostrophes, so we avoid it here).\n$target = $d.Range($searchRange.Start, $searchRange.End)\n$target.Text = $newText\n"}
